$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 527.98
    3  = 507.17
    4  = 493.17
    5  = 462.87
    6  = 415.53
    7  = 347.75
    8  = 324.88
    9  = 202.65
    10 = 164.88
    11 = 163.17
    12 = 144.22
    13 = 127.98
    14 = 87.33
    15 = 42.22
    16 = 22
    17 = 11.4
    18 = 5.55
    19 = 5.33
    20 = 4
    21 = 3.63
    22 = 3.33
    23 = 3.33
    24 = 3.33
    25 = 2.3
    26 = 2.3
    27 = 2
    28 = 2
    29 = 2
    30 = 2
    31 = 2
    32 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("K$row").Value = $values[$row]
}
